$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.242.85"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.657.98"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.36"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5319"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2629"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.50"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07846"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.537"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "1.638.80"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "1.885.40"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5510"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "0.0₅8181"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.51"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "26.217.99"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.621"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.87"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.029"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.73"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1221"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.230"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.471"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05788"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.283"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.600"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9537"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.820"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8519"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.69"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").Value = "1.043.46"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("D45").Value = "1.798.03"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.912"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05159"
$ws.Range("E51").Value = "  +0.14%  "
